# The commit swaps the deck's two theme parts: the slide-master theme
# (ppt/theme/theme1.xml, "Integral") ends up holding the stock "Office
# Theme" color palette that used to live in the notes-master theme
# (ppt/theme/theme2.xml), while theme2.xml ends up holding the palette
# that used to be in theme1.xml. The font scheme and format scheme
# (fills/lines/effects) are byte-identical between the two theme parts
# already, so the only observable difference is the 12-slot color
# scheme (clrScheme) used by the slides' theme (theme1.xml).
#
# The PowerPoint object model exposes that color scheme read/write via
# Slide.ThemeColorScheme (12 items, in clrMap order: dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink) which maps straight onto
# ppt/theme/theme1.xml's <a:clrScheme>. Re-pointing every slot to the
# "Office Theme" palette reproduces that swap.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function RGBOf([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index order matches a:clrScheme: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = RGBOf($officeThemeColors[$i - 1])
}
